$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $searchText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.Contains($searchText)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1. "...his fourteenth birthday that Keith Buchant decided that of..."
#    -> "...his eighteenth birthday that Keith Buchant decided; of..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " his fourteenth birthday that Keith Buchant decided that o",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " his eighteenth birthday that Keith Buchant decided; o", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "...not amongst the afore mentioned 'best friends', leading..."
#    -> "...not amongst his 'best friends', leading..."
# ---------------------------------------------------------------------------
$quoteOpen = [char]0x2018
$quoteClose = [char]0x2019
$d.Content.Find.Execute(
    "the afore mentioned " + $quoteOpen + "best friends" + $quoteClose,
    $true, $false, $false, $false, $false, $true, 1, $false,
    "his " + $quoteOpen + "best friends" + $quoteClose, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Remove one of the two blank paragraphs between
#    "...Well, Skylar, it's like this" he began." and "Skylar Clarke - Story 2"
# ---------------------------------------------------------------------------
$wellIdx = Find-ParagraphIndex $d "Well, Skylar, it"
$blank1 = $d.Paragraphs.Item($wellIdx + 1)
$blank2 = $d.Paragraphs.Item($wellIdx + 2)
if ($blank1.Range.Text.Trim() -eq "" -and $blank2.Range.Text.Trim() -eq "") {
    $blank1.Range.Delete() | Out-Null
}

# ---------------------------------------------------------------------------
# 4. Restructure the "Harrison Grey" area: add a new paragraph right after the
#    "Harrison Grey" heading paragraph that introduces a new possible
#    beginning for the Keith/Skylar story, and move the hidden "_GoBack"
#    bookmark from the end of the "Mackenzie Anson - Story 2" paragraph into
#    this new paragraph (right after "When Keith Buchant ").
# ---------------------------------------------------------------------------

# 4a. Delete the existing (hidden) _GoBack bookmark.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete() | Out-Null

# 4b. Insert the new paragraph immediately after the "Harrison Grey" heading,
#     i.e. immediately before the first of the two following blank paragraphs.
$harrisonIdx = Find-ParagraphIndex $d "Harrison Grey"
$insertBeforePara = $d.Paragraphs.Item($harrisonIdx + 1)
$insertPoint = $d.Range($insertBeforePara.Range.Start, $insertBeforePara.Range.Start)
$newParaText = "When Keith Buchant had recounted the story behind his latest set of broken glasses he made sure to mention two critical points. Firstly: the woman involved was mind-blowingly beautiful and secondly: he had failed, at that point in time, to get her number. "
$insertPoint.InsertBefore($newParaText + "`r")

# 4c. Re-add the "_GoBack" bookmark inside the new paragraph, right after
#     "When Keith Buchant ".
$newParaIdx = $harrisonIdx + 1
$newPara = $d.Paragraphs.Item($newParaIdx)
$markerRange = $newPara.Range.Duplicate()
$markerRange.Find.Execute("When Keith Buchant ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bookmarkRange = $d.Range($markerRange.End, $markerRange.End)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
